# Slide 5, shape "TextBox 25" (the callout reading "Interest Group makes
# changes and notifies community via tdwg-content"):
#   1. Nudge its position.
#   2. Collapse the three separate runs (the word "tdwg" had been split out,
#      presumably for a spell-check exception) back into a single run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(11)

# --- 1. Reposition -------------------------------------------------------
# EMU targets from the target OOXML: off x="6169300" y="523919".
# PowerPoint's Left/Top are in points (1 pt = 12700 EMU); the host stores
# them with limited (float32) precision, so a direct "emu/12700" assignment
# can land one EMU short after round-tripping. Nudge upward in tiny
# increments, re-reading the property each time, until the persisted value
# reaches the target.
$targetLeftEmu = 6169300
$pt = $targetLeftEmu / 12700.0
for ($i = 0; $i -lt 50; $i++) {
    $sh.Left = $pt
    $readBackEmu = [Math]::Round($sh.Left * 12700)
    if ($readBackEmu -ge $targetLeftEmu) { break }
    $pt += 0.000001
}

$targetTopEmu = 523919
$pt = $targetTopEmu / 12700.0
for ($i = 0; $i -lt 50; $i++) {
    $sh.Top = $pt
    $readBackEmu = [Math]::Round($sh.Top * 12700)
    if ($readBackEmu -ge $targetTopEmu) { break }
    $pt += 0.000001
}

# --- 2. Merge the text runs ----------------------------------------------
# Assigning the already-displayed text back to TextRange.Text is a no-op
# here (it already reads the same concatenated string), so first swap in a
# throwaway value to force a genuine replacement, then set the real text.
# That collapses "Interest Group makes changes and notifies community via "
# + "tdwg" + "-content" into one run, picking up the first run's formatting.
$sh.TextFrame.TextRange.Text = "x"
$sh.TextFrame.TextRange.Text = "Interest Group makes changes and notifies community via tdwg-content"
